$wb = $excel.ActiveWorkbook

# --- Sheet "snapshot": refresh scraped_at timestamps in column K (rows 2-37) ---
$snapshot = $wb.Worksheets.Item("snapshot")

$snapshot.Range("K2").Value = "2025-12-19T03:01:19.593792+00:00"
$snapshot.Range("K3").Value = "2025-12-19T03:01:21.526750+00:00"
$snapshot.Range("K4").Value = "2025-12-19T03:01:21.526785+00:00"
$snapshot.Range("K5").Value = "2025-12-19T03:01:21.526809+00:00"
$snapshot.Range("K6").Value = "2025-12-19T03:01:23.424599+00:00"
$snapshot.Range("K7").Value = "2025-12-19T03:01:25.782905+00:00"
$snapshot.Range("K8").Value = "2025-12-19T03:01:27.703548+00:00"
$snapshot.Range("K9").Value = "2025-12-19T03:01:27.703583+00:00"
$snapshot.Range("K10").Value = "2025-12-19T03:01:30.039079+00:00"
$snapshot.Range("K11").Value = "2025-12-19T03:01:34.758316+00:00"
$snapshot.Range("K12").Value = "2025-12-19T03:01:34.758345+00:00"
$snapshot.Range("K13").Value = "2025-12-19T03:01:37.166875+00:00"
$snapshot.Range("K14").Value = "2025-12-19T03:01:39.574468+00:00"
$snapshot.Range("K15").Value = "2025-12-19T03:01:41.846155+00:00"
$snapshot.Range("K16").Value = "2025-12-19T03:01:44.241245+00:00"
$snapshot.Range("K17").Value = "2025-12-19T03:01:44.241278+00:00"
$snapshot.Range("K18").Value = "2025-12-19T03:01:44.241297+00:00"
$snapshot.Range("K19").Value = "2025-12-19T03:01:44.241316+00:00"
$snapshot.Range("K20").Value = "2025-12-19T03:01:46.123505+00:00"
$snapshot.Range("K21").Value = "2025-12-19T03:01:46.123536+00:00"
$snapshot.Range("K22").Value = "2025-12-19T03:01:46.123558+00:00"
$snapshot.Range("K23").Value = "2025-12-19T03:01:48.057206+00:00"
$snapshot.Range("K24").Value = "2025-12-19T03:01:48.057236+00:00"
$snapshot.Range("K25").Value = "2025-12-19T03:01:48.057255+00:00"
$snapshot.Range("K26").Value = "2025-12-19T03:01:50.450516+00:00"
$snapshot.Range("K27").Value = "2025-12-19T03:01:50.450544+00:00"
$snapshot.Range("K28").Value = "2025-12-19T03:01:52.323874+00:00"
$snapshot.Range("K29").Value = "2025-12-19T03:01:52.323904+00:00"
$snapshot.Range("K30").Value = "2025-12-19T03:01:52.323922+00:00"
$snapshot.Range("K31").Value = "2025-12-19T03:01:54.665992+00:00"
$snapshot.Range("K32").Value = "2025-12-19T03:01:56.997677+00:00"
$snapshot.Range("K33").Value = "2025-12-19T03:01:56.997709+00:00"
$snapshot.Range("K34").Value = "2025-12-19T03:02:01.381420+00:00"
$snapshot.Range("K35").Value = "2025-12-19T03:02:01.381449+00:00"
$snapshot.Range("K36").Value = "2025-12-19T03:02:03.742081+00:00"
$snapshot.Range("K37").Value = "2025-12-19T03:02:03.742107+00:00"

# --- Sheet "new_injured": the single injury row (row 2) was processed/removed ---
$newInjured = $wb.Worksheets.Item("new_injured")
$newInjured.Rows.Item(2).Delete()

